# API 함수 개발#5 (revert of revert)
# - Uppercase the itemPart (column B) and itemGrade (column C) values used as
#   dropdown/category labels, fixing the "Sheld"/"Helmat" typos to "SHIELD"/"HELMET".
# - Rename header K1 from "itemSheldGager" to "itemSHIELDGager".
# - Update the two list data-validations (column B and C) to the new uppercase
#   option lists and extend them to the full column.
# - Widen column C slightly so the longer "itemGrade"/"LEGENDARY" text fits.
# - Restore the working zoom/selection state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header rename (column K, row 1) ---------------------------------
$ws.Cells.Item(1, 11).Value = "itemSHIELDGager"

# --- 2. Uppercase the Part (B) / Grade (C) columns for every data row ----
$partMap = @{
    "Weapon" = "WEAPON"
    "Gloves" = "GLOVES"
    "Shoes"  = "SHOES"
    "Sheld"  = "SHIELD"
    "Helmat" = "HELMET"
    "Armor"  = "ARMOR"
}
$gradeMap = @{
    "Normal"    = "NORMAL"
    "Rare"      = "RARE"
    "Epic"      = "EPIC"
    "Unique"    = "UNIQUE"
    "Legendary" = "LEGENDARY"
}

$lastRow = 91
for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = [string]$bCell.Value2
    if ($partMap.ContainsKey($bVal)) {
        $bCell.Value = $partMap[$bVal]
    }

    $cCell = $ws.Cells.Item($r, 3)
    $cVal = [string]$cCell.Value2
    if ($gradeMap.ContainsKey($cVal)) {
        $cCell.Value = $gradeMap[$cVal]
    }
}

# --- 3. Data validation lists (uppercased, extended to whole column) -----
$partRange = $ws.Range("B1:B1048576")
$partRange.Validation.Delete()
$partRange.Validation.Add(3, 1, 1, '"NONE, WEAPON, GLOVES, SHOES, SHIELD, HELMET, ARMOR, COUNT"')

$gradeRange = $ws.Range("C1:C1048576")
$gradeRange.Validation.Delete()
$gradeRange.Validation.Add(3, 1, 1, '"NONE, NORMAL, RARE, EPIC, UNIQUE, LEGENDARY, COUNT"')

# --- 4. Column width tweak ------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 12.43

# --- 5. Sheet view state (zoom + selection) -------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("F94").Select()
